# The "Förändrad" (Modified) date column (C) for every data row (2-216)
# was bumped from serial 45175 (2023-09-06) to serial 45177 (2023-09-08).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2:C216").Value = 45177
